$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column F (dSF) by row, per repulled/pushed data and mean calc
$values = @{
    2  = -3
    3  = 1
    4  = -1
    5  = 2
    6  = 8
    7  = 1
    8  = 2
    9  = 1
    10 = 5
    11 = 4
    12 = -1
    13 = 4
    14 = -1
    15 = 8
    16 = 1
    17 = -2
    18 = 6
    19 = -2
    20 = -4
    21 = 6
    22 = 4
    23 = -1
    24 = -2
    25 = 3
    26 = -1
    27 = -3
    29 = 1
    30 = 3
    31 = 0
    33 = -3
}

foreach ($row in $values.Keys) {
    $ws.Range("F$row").Value = $values[$row]
}
